$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange

# Replace distance codes first, then the size code, so no collisions occur.
# (Target tokens D55/D69/D86/S31 do not already exist in the source data.)
$usedRange.Replace("D51", "D55", -4123)
$usedRange.Replace("D64", "D69", -4123)
$usedRange.Replace("D80", "D86", -4123)
$usedRange.Replace("S30", "S31", -4123)
